# Auto-generated Excel COM-interop script to apply market-data refresh
# to the Halicarnassus Profits workbook (columns H-N per row).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 7186.5312

$ws.Range("H72").Value = 7186.5312

$ws.Range("H112").Value = 3237.7144
$ws.Range("J112").Value = 3237.7144
$ws.Range("L112").Value = 9713.143199999999
$ws.Range("N112").Value = -11929.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5319.606
$ws.Range("I32").Value = 5319.606
$ws.Range("K32").Value = 5319.606
$ws.Range("M32").Value = -5032.606

$ws.Range("H45").Value = 3252.7
$ws.Range("I45").Value = 2538.6667
$ws.Range("K45").Value = 2538.6667
$ws.Range("M45").Value = -2161.6667

$ws.Range("H110").Value = 2759.7778
$ws.Range("I110").Value = 1173.7222
$ws.Range("J110").Value = 4345.8335
$ws.Range("K110").Value = 1173.7222
$ws.Range("L110").Value = 4345.8335
$ws.Range("M110").Value = 871.2778000000001
$ws.Range("N110").Value = -8435.833500000001

$ws.Range("H132").Value = 5263.364
$ws.Range("I132").Value = 4210.778
$ws.Range("K132").Value = 12632.334
$ws.Range("M132").Value = -10102.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 13332.167
$ws.Range("I20").Value = 13332.167
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 13332.167
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -13085.167
$ws.Range("N20").ClearContents()

$ws.Range("H86").Value = 6217.4287
$ws.Range("I86").Value = 4875
$ws.Range("K86").Value = 4875
$ws.Range("M86").Value = -3752

$ws.Range("H89").Value = 6217.4287
$ws.Range("I89").Value = 4875
$ws.Range("K89").Value = 24375
$ws.Range("M89").Value = -18759

$ws.Range("H134").Value = 6298.1816
$ws.Range("J134").Value = 8293.6
$ws.Range("L134").Value = 24880.8
$ws.Range("N134").Value = -29950.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6304.515
$ws.Range("I31").Value = 2551.6667
$ws.Range("K31").Value = 2551.6667
$ws.Range("M31").Value = -2256.6667

$ws.Range("H34").Value = 6304.515
$ws.Range("I34").Value = 2551.6667
$ws.Range("K34").Value = 2551.6667
$ws.Range("M34").Value = -2349.6667

$ws.Range("H42").Value = 15000
$ws.Range("I42").Value = 15000
$ws.Range("K42").Value = 15000
$ws.Range("M42").Value = -14407

$ws.Range("H105").Value = 1573.8889
$ws.Range("J105").Value = 2722
$ws.Range("L105").Value = 2722
$ws.Range("N105").Value = -6216

$ws.Range("H134").Value = 3237.4285
$ws.Range("I134").Value = 1904.7
$ws.Range("K134").Value = 5714.1
$ws.Range("M134").Value = -3179.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 999.5
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 999.5
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 2998.5
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -3446.5

$ws.Range("H131").Value = 976.3333
$ws.Range("I131").Value = 976.3333
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 2928.9999
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 2111.0001
$ws.Range("N131").ClearContents()

$ws.Range("H140").Value = 2072.8823
$ws.Range("I140").Value = 1452.7142
$ws.Range("J140").Value = 4967
$ws.Range("K140").Value = 4358.142599999999
$ws.Range("L140").Value = 14901
$ws.Range("M140").Value = 821.8574000000008
$ws.Range("N140").Value = -25261

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 5833969.5
$ws.Range("I3").Value = 7778001
$ws.Range("J3").Value = 1875
$ws.Range("K3").Value = 7778001
$ws.Range("L3").Value = 1875
$ws.Range("M3").Value = -7777885
$ws.Range("N3").Value = -2107

$ws.Range("H7").Value = 6252375
$ws.Range("I7").Value = 10001333
$ws.Range("J7").Value = 4003000
$ws.Range("K7").Value = 10001333
$ws.Range("L7").Value = 4003000
$ws.Range("M7").Value = -10001221
$ws.Range("N7").Value = -4003224

$ws.Range("H8").Value = 6252375
$ws.Range("I8").Value = 10001333
$ws.Range("J8").Value = 4003000
$ws.Range("K8").Value = 10001333
$ws.Range("L8").Value = 4003000
$ws.Range("M8").Value = -10001194
$ws.Range("N8").Value = -4003278

$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()

$ws.Range("H11").Value = 6666836.5
$ws.Range("I11").Value = 6222422
$ws.Range("J11").Value = 8666700
$ws.Range("K11").Value = 6222422
$ws.Range("L11").Value = 8666700
$ws.Range("M11").Value = -6222283
$ws.Range("N11").Value = -8666978

$ws.Range("H20").Value = 35117.668
$ws.Range("I20").Value = 9950
$ws.Range("J20").Value = 38263.625
$ws.Range("K20").Value = 9950
$ws.Range("L20").Value = 38263.625
$ws.Range("M20").Value = -9705
$ws.Range("N20").Value = -38753.625

$ws.Range("H21").Value = 10333.667
$ws.Range("I21").Value = 7000.5
$ws.Range("J21").Value = 17000
$ws.Range("K21").Value = 7000.5
$ws.Range("L21").Value = 17000
$ws.Range("M21").Value = -6827.5
$ws.Range("N21").Value = -17346

$ws.Range("H30").Value = 10333.667
$ws.Range("I30").Value = 7000.5
$ws.Range("J30").Value = 17000
$ws.Range("K30").Value = 7000.5
$ws.Range("L30").Value = 17000
$ws.Range("M30").Value = -6895.5
$ws.Range("N30").Value = -17210

$ws.Range("H33").Value = 9999.5
$ws.Range("J33").Value = 9999.5
$ws.Range("L33").Value = 9999.5
$ws.Range("N33").Value = -10503.5

$ws.Range("H35").Value = 38017
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 38017
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 38017
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -38613

$ws.Range("H36").Value = 2666.5
$ws.Range("I36").Value = 1124.75
$ws.Range("J36").Value = 5750
$ws.Range("K36").Value = 1124.75
$ws.Range("L36").Value = 5750
$ws.Range("M36").Value = -639.75
$ws.Range("N36").Value = -6720

$ws.Range("H97").Value = 1051.1818
$ws.Range("I97").Value = 1007.55554
$ws.Range("K97").Value = 1007.55554
$ws.Range("M97").Value = -511.55554

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 9990
$ws.Range("I23").Value = 9990
$ws.Range("K23").Value = 9990
$ws.Range("M23").Value = -9760

$ws.Range("H33").Value = 9990
$ws.Range("I33").Value = 9990
$ws.Range("K33").Value = 9990
$ws.Range("M33").Value = -9700

$ws.Range("H42").Value = 9990
$ws.Range("I42").Value = 9990
$ws.Range("K42").Value = 9990
$ws.Range("M42").Value = -9427

$ws.Range("H43").Value = 6003
$ws.Range("I43").Value = 4006
$ws.Range("J43").Value = 8000
$ws.Range("K43").Value = 4006
$ws.Range("L43").Value = 8000
$ws.Range("M43").Value = -3813
$ws.Range("N43").Value = -8386

$ws.Range("H49").Value = 9990
$ws.Range("I49").Value = 9990
$ws.Range("K49").Value = 9990
$ws.Range("M49").Value = -9843

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 3241.3333
$ws.Range("I2").Value = 3241.3333
$ws.Range("K2").Value = 3241.3333
$ws.Range("M2").Value = -3129.3333

$ws.Range("H4").Value = 13855.875
$ws.Range("I4").Value = 20059.4
$ws.Range("K4").Value = 20059.4
$ws.Range("M4").Value = -19946.4

$ws.Range("H95").Value = 30714.143
$ws.Range("J95").Value = 30714.143
$ws.Range("L95").Value = 30714.143
$ws.Range("N95").Value = -36206.143
